{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the content of: async (context) => { ... }\n\nconst replacements = [\n  {\n    find: \"Play Big Win 777 Slot Game for Free - Review and Features\",\n    replace: \"Play Big Win 777 - Free Online Slot Game\"\n  },\n  {\n    find: \"Classic vintage graphics\",\n    replace: \"Vintage graphics and Las Vegas atmosphere\"\n  },\n  {\n    find: \"Chance Wheel feature with big payouts\",\n    replace: \"Classic symbols like BAR and colorful 7s\"\n  },\n  {\n    find: \"Big payouts with multipliers\",\n    replace: \"Chance Wheel feature with big payouts up to 777x\"\n  },\n  {\n    find: \"Similar slot games available online\",\n    replace: \"Similar slot games available for more options\"\n  },\n  {\n    find: \"Lacks advanced features like free spins and bonus rounds\",\n    replace: \"Limited variety of symbols\"\n  },\n  {\n    find: \"Limited paylines\",\n    replace: \"No free spins feature\"\n  },\n  {\n    find: \"Read our detailed review of Big Win 777 online slot game. Play for free and learn about the vintage graphics, Chance Wheel feature, and big payouts.\",\n    replace: \"Detailed review of Big Win 777, a traditional slot game with vintage graphics and big payouts. Play for free.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Big Win 777 Slot Game for Free - Review and Features\"; Replace = \"Play Big Win 777 - Free Online Slot Game\" },\n    @{ Find = \"Classic vintage graphics\"; Replace = \"Vintage graphics and Las Vegas atmosphere\" },\n    @{ Find = \"Chance Wheel feature with big payouts\"; Replace = \"Classic symbols like BAR and colorful 7s\" },\n    @{ Find = \"Big payouts with multipliers\"; Replace = \"Chance Wheel feature with big payouts up to 777x\" },\n    @{ Find = \"Similar slot games available online\"; Replace = \"Similar slot games available for more options\" },\n    @{ Find = \"Lacks advanced features like free spins and bonus rounds\"; Replace = \"Limited variety of symbols\" },\n    @{ Find = \"Limited paylines\"; Replace = \"No free spins feature\" },\n    @{ Find = \"Read our detailed review of Big Win 777 online slot game. Play for free and learn about the vintage graphics, Chance Wheel feature, and big payouts.\"; Replace = \"Detailed review of Big Win 777, a traditional slot game with vintage graphics and big payouts. Play for free.\" }\n)\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $item.Find\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $item.Replace\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n$d.Save()\n"}
